# Replace the pandoc-style title block (italic subtitle + "Chapter 2..." line,
# then a bold "By Dorothy Day" byline paragraph) with a single plain
# "% Dorothy Day" line, as used by pandoc title blocks.
#
# Before:
#   P1: "From Union Square to Rome" (italic) + "," + " " + <br/> +
#       "Chapter 2 - Childhood =============================="
#   P2: "By Dorothy Day" (bold)
#   P3: "1938, Chapter 2, pp. 18-27."
#
# After:
#   P1: "% Dorothy Day"   (plain run, no formatting)
#   P2: "1938, Chapter 2, pp. 18-27."   (unchanged)

$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)

# Sanity-check we're targeting the right content before mutating anything.
$p1Text = $p1.Range.Text
$p2Text = $p2.Range.Text
if ($p1Text -notmatch "From Union Square to Rome" -or $p2Text -notmatch "By Dorothy Day") {
    throw "Unexpected document content; aborting to avoid corrupting the document."
}

# Range spanning both paragraphs (including both paragraph marks), to be
# replaced wholesale by a single new paragraph.
$targetRange = $d.Range($p1.Range.Start, $p2.Range.End)

$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                    '<w:body>' +
                        '<w:p>' +
                            '<w:r>' +
                                '<w:t xml:space="preserve">% Dorothy Day</w:t>' +
                            '</w:r>' +
                        '</w:p>' +
                    '</w:body>' +
                '</w:document>' +
            '</pkg:xmlData>' +
        '</pkg:part>' +
    '</pkg:package>'

$targetRange.InsertXML($newXml)
